$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @(17.93638541084523,6.267222904724977,3.776792924455576,10.49183963931413,56.73234970993249,7.344005520526261,41.13147533566484,10.31861299239337,16.29887471242875,11.88618270867794)
  3 = @(17.89379833172608,6.14764753714271,3.795759176706595,10.51373394722424,56.54239607837081,7.344005520526261,41.03273893390486,10.33078947304563,16.27115215239628,11.9054666500998)
  4 = @(17.87280340412057,6.075128913113524,3.808166578393005,10.52828887710655,56.43286028250721,7.344005520526261,40.97632699897966,10.33886500076165,16.25841551737983,11.91909959623421)
  5 = @(17.86555192232207,6.045848442109947,3.81341485492007,10.53450023911106,56.39002711416634,7.344005520526261,40.95440529428479,10.34230680158095,16.25430763588216,11.925106087725)
  6 = @(17.8644267958335,6.041004148197964,3.81429794905819,10.53554856526464,56.38302415199469,7.344005520526261,40.95082983117435,10.34288743639969,16.25369102236762,11.92613070279146)
  7 = @(17.87270031805618,6.074732868051539,3.808236579714392,10.52837151081189,56.43227529315862,7.344005520526261,40.97602702754204,10.3389108064659,16.2583557287895,11.91917877574251)
  8 = @(17.92063707169208,6.225833167128676,3.783174647639835,10.49915838940626,56.66539277675368,7.344005520526261,41.09656130428847,10.32268727848368,16.28842981370644,11.89245991162393)
  9 = @(18.05510837660585,6.527333364683563,3.740050810041367,10.45066811854134,57.17793648595563,7.344005520526261,41.36600640648423,10.2956132673555,16.38114755841291,11.85427695545831)
  10 = @(18.17791207231304,6.749443734800667,3.712006463408541,10.42037164457594,57.58696627340909,7.344005520526261,41.58363032356536,10.27859354224366,16.46942855084463,11.83487293966066)
  11 = @(18.23882206970784,6.850098566239938,3.700031606282018,10.40773914105895,57.77978099216221,7.344005520526261,41.68678013741286,10.27147050443684,16.51386016131485,11.82791925574073)
  12 = @(18.26259641209458,6.888117779970725,3.69560905683187,10.40312027599728,57.85373559035412,7.344005520526261,41.72642542281064,10.26886194863951,16.53128899963651,11.82555497256776)
  13 = @(18.25744492333568,6.879934593807134,3.696556555198408,10.40410771026645,57.83776683859997,7.344005520526261,41.71786130067426,10.26941980345442,16.52750872552877,11.8260522105817)
  14 = @(18.24076388788515,6.853228594962382,3.699665517149316,10.40735584453216,57.78584662291711,7.344005520526261,41.6900301339307,10.2712541191803,16.51528201120035,11.82771935788262)
  15 = @(18.23063809908773,6.83685662966175,3.701584426796872,10.40936686595324,57.75416547029055,7.344005520526261,41.67305850685489,10.27238924472871,16.50787106057384,11.82877554132679)
  16 = @(18.17403161112664,6.742854393032281,3.712804756700199,10.42122029606322,57.57449856753746,7.344005520526261,41.57697174631058,10.27907148842493,16.46660990660884,11.83536503765174)
  17 = @(18.14058636894462,6.685058600093766,3.719888182085376,10.42878604770298,57.46598537277278,7.344005520526261,41.51908043435187,10.2833292563951,16.44238434666198,11.83988697583986)
  18 = @(18.12182557264988,6.651782167769352,3.724036076451977,10.43324589615144,57.40420822240858,7.344005520526261,41.48617443532478,10.28583652159873,16.4288532641467,11.84266423009652)
  19 = @(18.11555573114327,6.640510822995759,3.725453154082797,10.43477452799804,57.38340179867758,7.344005520526261,41.47510064180406,10.2866954616042,16.42434136998077,11.84363486252527)
  20 = @(18.14409752174261,6.691214867076198,3.719126515630369,10.4279694630107,57.47747109789157,7.344005520526261,41.52520262614345,10.28286997689002,16.44492158391613,11.83938735882669)
  21 = @(18.24564440688559,6.86107573465448,3.698749301902536,10.40639732047082,57.80107158845854,7.344005520526261,41.69818905560672,10.27071292878881,16.51885700372517,11.82722238165686)
  22 = @(18.31613553322013,6.971508969194848,3.686084622665763,10.39325896328184,58.01802642340142,7.344005520526261,41.81464749031741,10.26328497015265,16.5706904089223,11.82083911696213)
  23 = @(18.2781414886687,6.9126349174519,3.692784405652426,10.40018345386436,57.9017439662886,7.344005520526261,41.75218441876397,10.26720216172273,16.54270835793976,11.82410274798806)
  24 = @(18.14250867268111,6.688431767806175,3.719470629714055,10.42833829766885,57.47227650372916,7.344005520526261,41.52243360954165,10.28307743205622,16.44377326305304,11.83961268266006)
  25 = @(18.01446088613126,6.445483101346332,3.75107560713835,10.46284772993058,57.03347723196368,7.344005520526261,41.28962188534969,10.30243191578427,16.35249087744282,11.8630862563232)
}

$cols = @(2,3,4,5,6,8,9,10,11,12)  # B,C,D,E,F,H,I,J,K,L
foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Cells.Item($r, $cols[$i]).Value = $vals[$i]
  }
}

$wb.Save()